$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Rename the sheet: "RO & CO Hearing Allocation" -> "RO Allocations"
$ws.Name = "RO Allocations"

# Update the title cell to drop the Central Office mention
$ws.Range("A1").Value = "Allocation of Regional Office Video Hearings"

# Remove the "Central Office" allocation row (row 4); this shifts all
# subsequent rows up by one and shrinks the used range from N100 to N99
$ws.Rows.Item(4).Delete()
